$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2105263157894737
$ws.Range("C2").Value = 0.5131578947368421
$ws.Range("J2").Value = 0.009868421052631578
$ws.Range("P2").Value = 0.1710526315789474
$ws.Range("S2").Value = 0.09539473684210527
$ws.Range("B3").Value = 0.01176470588235294
$ws.Range("C3").Value = 0.02941176470588235
$ws.Range("J3").Value = 0.02352941176470588
$ws.Range("P3").Value = 0.7588235294117647
$ws.Range("S3").Value = 0.1764705882352941
$ws.Range("B6").Value = 0.02766798418972332
$ws.Range("D6").Value = 0.01185770750988142
$ws.Range("E6").Value = 0.003952569169960474
$ws.Range("F6").Value = 0.06324110671936758
$ws.Range("J6").Value = 0.2490118577075099
$ws.Range("O6").Value = 0.007905138339920948
$ws.Range("Q6").Value = 0.1897233201581028
$ws.Range("R6").Value = 0.09090909090909091
$ws.Range("S6").Value = 0.3557312252964427
$ws.Range("B7").Value = 0.096045197740113
$ws.Range("D7").Value = 0.02259887005649718
$ws.Range("E7").Value = 0.005649717514124294
$ws.Range("F7").Value = 0.05084745762711865
$ws.Range("J7").Value = 0.1186440677966102
$ws.Range("O7").Value = 0.005649717514124294
$ws.Range("Q7").Value = 0.1807909604519774
$ws.Range("R7").Value = 0.07344632768361582
$ws.Range("S7").Value = 0.4463276836158192
$ws.Range("B8").Value = 0.08921161825726141
$ws.Range("D8").Value = 0.01867219917012448
$ws.Range("E8").Value = 0.002074688796680498
$ws.Range("F8").Value = 0.05809128630705394
$ws.Range("J8").Value = 0.07883817427385892
$ws.Range("O8").Value = 0.01037344398340249
$ws.Range("Q8").Value = 0.1867219917012448
$ws.Range("R8").Value = 0.08713692946058091
$ws.Range("S8").Value = 0.4688796680497925
$ws.Range("B9").Value = 0.08928571428571429
$ws.Range("D9").Value = 0.025
$ws.Range("E9").Value = 0.003571428571428571
$ws.Range("F9").Value = 0.05357142857142857
$ws.Range("J9").Value = 0.075
$ws.Range("R9").Value = 0.06428571428571428
$ws.Range("S9").Value = 0.4892857142857143
$ws.Range("B10").Value = 0.1020114942528736
$ws.Range("D10").Value = 0.01293103448275862
$ws.Range("E10").Value = 0.0007183908045977011
$ws.Range("F10").Value = 0.07112068965517242
$ws.Range("J10").Value = 0.09985632183908046
$ws.Range("O10").Value = 0.01939655172413793
$ws.Range("Q10").Value = 0.2413793103448276
$ws.Range("R10").Value = 0.07758620689655173
$ws.Range("S10").Value = 0.375
$ws.Range("F11").Value = 0.00353356890459364
$ws.Range("G11").Value = 0.1201413427561837
$ws.Range("J11").Value = 0.1130742049469965
$ws.Range("K11").Value = 0.1908127208480565
$ws.Range("L11").Value = 0.5653710247349824
$ws.Range("S11").Value = 0.007067137809187279
$ws.Range("G12").Value = 0.7470588235294118
$ws.Range("J12").Value = 0.1588235294117647
$ws.Range("K12").Value = 0.02352941176470588
$ws.Range("L12").Value = 0.05294117647058823
$ws.Range("S12").Value = 0.01764705882352941
$ws.Range("G13").Value = 0.7941176470588235
$ws.Range("J13").Value = 0.1764705882352941
$ws.Range("S13").Value = 0.02941176470588235
$ws.Range("F15").Value = 0.02109704641350211
$ws.Range("H15").Value = 0.160337552742616
$ws.Range("I15").Value = 0.08438818565400844
$ws.Range("J15").Value = 0.3670886075949367
$ws.Range("K15").Value = 0.04641350210970464
$ws.Range("M15").Value = 0.008438818565400843
$ws.Range("O15").Value = 0.04641350210970464
$ws.Range("S15").Value = 0.2658227848101266
$ws.Range("F16").Value = 0.03980099502487562
$ws.Range("H16").Value = 0.1492537313432836
$ws.Range("I16").Value = 0.08955223880597014
$ws.Range("J16").Value = 0.472636815920398
$ws.Range("K16").Value = 0.0845771144278607
$ws.Range("M16").Value = 0.01990049751243781
$ws.Range("N16").Value = 0.004975124378109453
$ws.Range("O16").Value = 0.05472636815920398
$ws.Range("S16").Value = 0.0845771144278607
$ws.Range("F17").Value = 0.02142857142857143
$ws.Range("H17").Value = 0.1660714285714286
$ws.Range("I17").Value = 0.1125
$ws.Range("J17").Value = 0.4446428571428571
$ws.Range("K17").Value = 0.07142857142857142
$ws.Range("M17").Value = 0.0125
$ws.Range("O17").Value = 0.075
$ws.Range("S17").Value = 0.09642857142857143
$ws.Range("F18").Value = 0.04477611940298507
$ws.Range("H18").Value = 0.1393034825870647
$ws.Range("I18").Value = 0.1144278606965174
$ws.Range("J18").Value = 0.417910447761194
$ws.Range("K18").Value = 0.1144278606965174
$ws.Range("M18").Value = 0.01990049751243781
$ws.Range("N18").Value = 0.004975124378109453
$ws.Range("O18").Value = 0.05472636815920398
$ws.Range("S18").Value = 0.08955223880597014
$ws.Range("F19").Value = 0.0148936170212766
$ws.Range("H19").Value = 0.2099290780141844
$ws.Range("I19").Value = 0.1106382978723404
$ws.Range("J19").Value = 0.3858156028368794
$ws.Range("K19").Value = 0.09219858156028368
$ws.Range("M19").Value = 0.01347517730496454
$ws.Range("N19").Value = 0.0007092198581560284
$ws.Range("O19").Value = 0.07021276595744681
$ws.Range("S19").Value = 0.1021276595744681
